$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "alignment" setting was renamed to "indexing" (and its description updated
# to reflect that the bowtie index is generated, not just checked/performed).
$ws.Range("A7").Value = "indexing"
$ws.Range("B7").Value = "TRUE if bowtie index needs to be generated, else FALSE"

# Match the author's final selection (cell B7) on the active sheet.
[void]$ws.Range("B7").Select()
